# Fill in the "Temps" (time invested) and "Avancement" (progress) columns
# for the first four tasks of the "Sprint 1 - Bilan" worksheet, then make
# that sheet the active one (matching the author finishing off the
# "Rédaction du document de conception" task and leaving the workbook with
# the Bilan sheet selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1 - Bilan")

# Row 8  - "Valider le projet avec l'enseignant"   -> 0:15, 100%
$ws.Range("E8").Value = 0.010416666666666666
$ws.Range("F8").Value = 1

# Row 9  - "Réaliser les éléments de rédaction technique" -> 6:00, 100%
$ws.Range("E9").Value = 0.25
$ws.Range("F9").Value = 1

# Row 10 - "Rédaction du document de conception" -> 5:00, 100%
$ws.Range("E10").Value = 0.20833333333333301
$ws.Range("F10").Value = 1

# Row 11 - "Rédaction du document de planification" -> 0:45, 100%
$ws.Range("E11").Value = 0.03125
$ws.Range("F11").Value = 1

# Update the selection left on this sheet.
$ws.Range("F19").Select()

# Make "Sprint 1 - Bilan" the active sheet/tab (was "Sprint 1 - Planification").
$ws.Activate()

$wb.Application.Calculate()
